$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GDNN T11-T12")

# --- Fill in row 18 (ngay 5/12/2023) raw values ---
$ws.Range("B18").Value = 234284.0
$ws.Range("C18").Value = 3925541.0
$ws.Range("D18").Value = -3691257.0
$ws.Range("E18").Value = -83351620000.0
$ws.Range("F18").Value = 27463600.0
$ws.Range("J18").Value = 580500.0
$ws.Range("K18").Value = 2602174.0
$ws.Range("L18").Value = -2021674.0
$ws.Range("M18").Value = -67407980000.0
$ws.Range("N18").Value = 23475502.0
$ws.Range("R18").Value = 60000.0
$ws.Range("S18").Value = 1234086.0
$ws.Range("T18").Value = -1174086.0
$ws.Range("U18").Value = -104177060000.0
$ws.Range("V18").Value = 2339861.0
$ws.Range("Z18").Value = 1281300.0
$ws.Range("AA18").Value = 5233690.0
$ws.Range("AB18").Value = -3952390.0
$ws.Range("AC18").Value = -76291340000.0
$ws.Range("AD18").Value = 17250007.0
$ws.Range("AH18").Value = 151200.0
$ws.Range("AI18").Value = 1398700.0
$ws.Range("AJ18").Value = -1247500.0
$ws.Range("AK18").Value = -35020850000.0
$ws.Range("AL18").Value = 13499600.0
$ws.Range("AP18").Value = 614100.0
$ws.Range("AQ18").Value = 7498060.0
$ws.Range("AR18").Value = -6883960.0
$ws.Range("AS18").Value = -188056260000.0
$ws.Range("AT18").Value = 22656100.0
$ws.Range("AX18").Value = 135700.0
$ws.Range("AY18").Value = 129600.0
$ws.Range("AZ18").Value = 6100.0
$ws.Range("BA18").Value = 159870000.0
$ws.Range("BB18").Value = 10469500.0
$ws.Range("BF18").Value = 677700.0
$ws.Range("BG18").Value = 675878.0
$ws.Range("BH18").Value = 1822.0
$ws.Range("BI18").Value = 31370000.0
$ws.Range("BJ18").Value = 7196600.0
$ws.Range("BN18").Value = 1816750.0
$ws.Range("BO18").Value = 2133505.0
$ws.Range("BP18").Value = -316755.0
$ws.Range("BQ18").Value = -7268400000.0
$ws.Range("BR18").Value = 3961100.0
$ws.Range("BV18").Value = 111900.0
$ws.Range("BW18").Value = 234902.0
$ws.Range("BX18").Value = -123002.0
$ws.Range("BY18").Value = -12566610000.0
$ws.Range("BZ18").Value = 414200.0
$ws.Range("CD18").Value = 2325014.0
$ws.Range("CE18").Value = 6547281.0
$ws.Range("CF18").Value = -4222267.0
$ws.Range("CG18").Value = -172274840000.0
$ws.Range("CH18").Value = 10227900.0
$ws.Range("CL18").Value = 253500.0
$ws.Range("CM18").Value = 3893004.0
$ws.Range("CN18").Value = -3639504.0
$ws.Range("CO18").Value = -76270090000.0
$ws.Range("CP18").Value = 35007200.0
$ws.Range("CT18").Value = 172183.0
$ws.Range("CU18").Value = 219300.0
$ws.Range("CV18").Value = -47117.0
$ws.Range("CW18").Value = -1322000000.0
$ws.Range("CX18").Value = 10358600.0
$ws.Range("DB18").Value = 659980.0
$ws.Range("DC18").Value = 501401.0
$ws.Range("DD18").Value = 158579.0
$ws.Range("DE18").Value = 2941910000.0
$ws.Range("DF18").Value = 39037500.0
$ws.Range("DJ18").Value = 635900.0
$ws.Range("DK18").Value = 281700.0
$ws.Range("DL18").Value = 354200.0
$ws.Range("DM18").Value = 11599390000.0
$ws.Range("DN18").Value = 5979300.0
$ws.Range("DR18").Value = 62500.0
$ws.Range("DS18").Value = 155800.0
$ws.Range("DT18").Value = -93300.0
$ws.Range("DU18").Value = -6625120000.0
$ws.Range("DV18").Value = 948300.0
$ws.Range("DZ18").Value = 2300.0
$ws.Range("EA18").Value = 10090.0
$ws.Range("EB18").Value = -7790.0
$ws.Range("EC18").Value = -288210000.0
$ws.Range("ED18").Value = 3180500.0
$ws.Range("EH18").Value = 7500.0
$ws.Range("EI18").Value = 41200.0
$ws.Range("EJ18").Value = -33700.0
$ws.Range("EK18").Value = -2186870000.0
$ws.Range("EL18").Value = 1137000.0
$ws.Range("EP18").Value = 90200.0
$ws.Range("EQ18").Value = 258000.0
$ws.Range("ER18").Value = -167800.0
$ws.Range("ES18").Value = -4101900000.0
$ws.Range("ET18").Value = 11721200.0
$ws.Range("EX18").Value = 157400.0
$ws.Range("EY18").Value = 383421.0
$ws.Range("EZ18").Value = -226021.0
$ws.Range("FA18").Value = -5319820000.0
$ws.Range("FB18").Value = 28094000.0
$ws.Range("FF18").Value = 1600.0
$ws.Range("FG18").Value = 145985.0
$ws.Range("FH18").Value = -144385.0
$ws.Range("FI18").Value = -14115170000.0
$ws.Range("FJ18").Value = 2049000.0
$ws.Range("FN18").Value = 233200.0
$ws.Range("FO18").Value = 168400.0
$ws.Range("FP18").Value = 64800.0
$ws.Range("FQ18").Value = 4494730000.0
$ws.Range("FR18").Value = 1425500.0

$ws.Range("G18").Value = 0.151466850667793
$ws.Range("H18").Value = -0.006651884700665125
$ws.Range("O18").Value = 0.13557426801778297
$ws.Range("P18").Value = -0.010510510510510555
$ws.Range("W18").Value = 0.5530610579004479
$ws.Range("X18").Value = -0.007025761124121713
$ws.Range("AE18").Value = 0.3776804264485226
$ws.Range("AF18").Value = -0.018276762402088847
$ws.Range("AM18").Value = 0.11481080920916176
$ws.Range("AN18").Value = -0.01258992805755388
$ws.Range("AU18").Value = 0.3580563292005244
$ws.Range("AV18").Value = -0.009174311926605505
$ws.Range("BC18").Value = 0.025340274129614596
$ws.Range("BD18").Value = -0.016985138004246225
$ws.Range("BK18").Value = 0.18808576272128505
$ws.Range("BL18").Value = 0.0012515644555693905
$ws.Range("BS18").Value = 0.997262124157431
$ws.Range("BT18").Value = -0.006564551422319412
$ws.Range("CA18").Value = 0.83728150651859
$ws.Range("CB18").Value = -0.02161100196463657
$ws.Range("CI18").Value = 0.8674600846703624
$ws.Range("CJ18").Value = -0.024691358024691357
$ws.Range("CQ18").Value = 0.11844717658081766
$ws.Range("CR18").Value = 0.02179176755447956
$ws.Range("CY18").Value = 0.03779304153070878
$ws.Range("CZ18").Value = -0.01610017889087654
$ws.Range("DG18").Value = 0.02975039385206532
$ws.Range("DH18").Value = -0.01098901098901095
$ws.Range("DO18").Value = 0.15346277992407137
$ws.Range("DP18").Value = 0.0
$ws.Range("DW18").Value = 0.2302014130549404
$ws.Range("DX18").Value = -0.011396011396011355
$ws.Range("EE18").Value = 0.003895613897185977
$ws.Range("EF18").Value = -0.012178619756427488
$ws.Range("EM18").Value = 0.04283201407211961
$ws.Range("EN18").Value = -0.017080745341614818
$ws.Range("EU18").Value = 0.02970685595331536
$ws.Range("EV18").Value = 0.004098360655737763
$ws.Range("FC18").Value = 0.019250409340072612
$ws.Range("FD18").Value = 0.006410256410256502
$ws.Range("FK18").Value = 0.07202781844802343
$ws.Range("FL18").Value = -0.004110996916752371
$ws.Range("FS18").Value = 0.28172571027709575
$ws.Range("FT18").Value = 0.03047895500725681

# --- Apply percentage number format (0.00%) to GDNN/Total-KL and Bien dong(%) columns ---
$ws.Range("G18:H18").NumberFormat = "0.00%"
$ws.Range("O18:P18").NumberFormat = "0.00%"
$ws.Range("W18:X18").NumberFormat = "0.00%"
$ws.Range("AE18:AF18").NumberFormat = "0.00%"
$ws.Range("AM18:AN18").NumberFormat = "0.00%"
$ws.Range("AU18:AV18").NumberFormat = "0.00%"
$ws.Range("BC18:BD18").NumberFormat = "0.00%"
$ws.Range("BK18:BL18").NumberFormat = "0.00%"
$ws.Range("BS18:BT18").NumberFormat = "0.00%"
$ws.Range("CA18:CB18").NumberFormat = "0.00%"
$ws.Range("CI18:CJ18").NumberFormat = "0.00%"
$ws.Range("CQ18:CR18").NumberFormat = "0.00%"
$ws.Range("CY18:CZ18").NumberFormat = "0.00%"
$ws.Range("DG18:DH18").NumberFormat = "0.00%"
$ws.Range("DO18:DP18").NumberFormat = "0.00%"
$ws.Range("DW18:DX18").NumberFormat = "0.00%"
$ws.Range("EE18:EF18").NumberFormat = "0.00%"
$ws.Range("EM18:EN18").NumberFormat = "0.00%"
$ws.Range("EU18:EV18").NumberFormat = "0.00%"
$ws.Range("FC18:FD18").NumberFormat = "0.00%"
$ws.Range("FK18:FL18").NumberFormat = "0.00%"
$ws.Range("FS18:FT18").NumberFormat = "0.00%"

# --- Restore frozen-pane scroll position to column B (was parked at EW1) ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2

